$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the original sheet from "Sheet1" to "data"
$ws.Name = "data"

# Add a new worksheet right after "data" and name it "data_dummy"
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws)
$ws2.Name = "data_dummy"

# Keep "data" as the active/selected sheet (matches tabSelected="1" on sheet1 in the diff)
$ws.Activate()

Write-Host "Sheets:" $wb.Worksheets.Count
